$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 5418.6
$ws.Range("I18").Value = 2800
$ws.Range("J18").Value = 7164.3335
$ws.Range("K18").Value = 2800
$ws.Range("L18").Value = 7164.3335
$ws.Range("M18").Value = -2516
$ws.Range("N18").Value = -7732.3335
$ws.Range("H86").Value = 7245.091
$ws.Range("I86").Value = 6310.5
$ws.Range("J86").Value = 8366.6
$ws.Range("K86").Value = 6310.5
$ws.Range("L86").Value = 8366.6
$ws.Range("M86").Value = -5187.5
$ws.Range("N86").Value = -10612.6
$ws.Range("H89").Value = 7245.091
$ws.Range("I89").Value = 6310.5
$ws.Range("J89").Value = 8366.6
$ws.Range("K89").Value = 31552.5
$ws.Range("L89").Value = 41833
$ws.Range("M89").Value = -25936.5
$ws.Range("N89").Value = -53065
$ws.Range("H137").Value = 1860.9333
$ws.Range("I137").Value = 1676.1578
$ws.Range("J137").Value = 2180.0908
$ws.Range("K137").Value = 5028.4734
$ws.Range("L137").Value = 6540.2724
$ws.Range("M137").Value = -2478.4734
$ws.Range("N137").Value = -11640.2724
$ws.Range("H138").Value = 1856307.2
$ws.Range("J138").Value = 2852317.8
$ws.Range("L138").Value = 8556953.399999999
$ws.Range("N138").Value = -8567233.399999999
$ws.Range("H141").Value = 4518.05
$ws.Range("I141").Value = 2478.6365
$ws.Range("K141").Value = 7435.9095
$ws.Range("M141").Value = -2255.9095

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3204.0344
$ws.Range("I2").Value = 3024.95
$ws.Range("K2").Value = 3024.95
$ws.Range("M2").Value = -2911.95
$ws.Range("H37").Value = 49997.5
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = 49997.5
$ws.Range("K37").Value = 0
$ws.Range("L37").Value = 49997.5
$ws.Range("M37").ClearContents()
$ws.Range("N37").Value = -50543.5
$ws.Range("H74").Value = 2397.242
$ws.Range("I74").Value = 1967.95
$ws.Range("K74").Value = 1967.95
$ws.Range("M74").Value = -1093.95
$ws.Range("H77").Value = 2397.242
$ws.Range("I77").Value = 1967.95
$ws.Range("K77").Value = 9839.75
$ws.Range("M77").Value = -5471.75
$ws.Range("H110").Value = 957.4286
$ws.Range("I110").Value = 838.8077
$ws.Range("K110").Value = 838.8077
$ws.Range("M110").Value = 1206.1923
$ws.Range("H116").Value = 3204.0344
$ws.Range("I116").Value = 3024.95
$ws.Range("K116").Value = 3024.95
$ws.Range("M116").Value = -730.9499999999998
$ws.Range("H122").Value = 10104755
$ws.Range("I122").Value = 18522234
$ws.Range("K122").Value = 55566702
$ws.Range("M122").Value = -55564252
$ws.Range("H132").Value = 2352.6128
$ws.Range("I132").Value = 2181.9614
$ws.Range("K132").Value = 6545.8842
$ws.Range("M132").Value = -4015.8842

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3204.0344
$ws.Range("I3").Value = 3024.95
$ws.Range("K3").Value = 3024.95
$ws.Range("M3").Value = -2910.95
$ws.Range("H134").Value = 3089.25
$ws.Range("I134").Value = 3570.6
$ws.Range("J134").Value = 2287
$ws.Range("K134").Value = 10711.8
$ws.Range("L134").Value = 6861
$ws.Range("M134").Value = -8176.799999999999
$ws.Range("N134").Value = -11931

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 3925.389
$ws.Range("I105").Value = 6751.5
$ws.Range("K105").Value = 6751.5
$ws.Range("M105").Value = -5004.5
$ws.Range("H107").Value = 844.55554
$ws.Range("I107").Value = 1037.1428
$ws.Range("J107").Value = 722
$ws.Range("K107").Value = 1037.1428
$ws.Range("L107").Value = 722
$ws.Range("M107").Value = 882.8571999999999
$ws.Range("N107").Value = -4562
$ws.Range("H132").Value = 1286.4
$ws.Range("I132").Value = 1303.2858
$ws.Range("J132").Value = 1050
$ws.Range("K132").Value = 3909.8574
$ws.Range("L132").Value = 3150
$ws.Range("M132").Value = -1379.8574
$ws.Range("N132").Value = -8210

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H36").Value = 754.2
$ws.Range("I36").Value = 754.2
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 2262.6
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = -2093.6
$ws.Range("N36").ClearContents()
$ws.Range("H87").Value = 0
$ws.Range("I87").Value = 0
$ws.Range("K87").Value = 0
$ws.Range("M87").ClearContents()
$ws.Range("H90").Value = 0
$ws.Range("I90").Value = 0
$ws.Range("K90").Value = 0
$ws.Range("M90").ClearContents()
$ws.Range("H97").Value = 239
$ws.Range("I97").Value = 203.33333
$ws.Range("K97").Value = 609.99999
$ws.Range("M97").Value = -113.99999
$ws.Range("H103").Value = 44
$ws.Range("I103").Value = 46
$ws.Range("J103").Value = 41
$ws.Range("K103").Value = 138
$ws.Range("L103").Value = 123
$ws.Range("M103").Value = 741
$ws.Range("N103").Value = -1881
$ws.Range("H132").Value = 2617.5293
$ws.Range("I132").Value = 1725.25
$ws.Range("J132").Value = 2892.077
$ws.Range("K132").Value = 15527.25
$ws.Range("L132").Value = 26028.693
$ws.Range("M132").Value = -12997.25
$ws.Range("N132").Value = -31088.693
$ws.Range("H139").Value = 4670.5386
$ws.Range("I139").Value = 4422.8184
$ws.Range("K139").Value = 13268.4552
$ws.Range("M139").Value = -8128.4552

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 22225638
$ws.Range("I70").Value = 37039616
$ws.Range("J70").Value = 4668.1665
$ws.Range("K70").Value = 37039616
$ws.Range("L70").Value = 4668.1665
$ws.Range("M70").Value = -37039346
$ws.Range("N70").Value = -5208.1665
$ws.Range("H73").Value = 22225638
$ws.Range("I73").Value = 37039616
$ws.Range("J73").Value = 4668.1665
$ws.Range("K73").Value = 37039616
$ws.Range("L73").Value = 4668.1665
$ws.Range("M73").Value = -37038680
$ws.Range("N73").Value = -6540.1665
$ws.Range("H113").Value = 7619.357
$ws.Range("I113").Value = 7566.6313
$ws.Range("K113").Value = 7566.6313
$ws.Range("M113").Value = -5396.6313
$ws.Range("H122").Value = 3956.3235
$ws.Range("I122").Value = 2608.9583
$ws.Range("J122").Value = 7190
$ws.Range("K122").Value = 7826.874899999999
$ws.Range("L122").Value = 21570
$ws.Range("M122").Value = -5376.874899999999
$ws.Range("N122").Value = -26470
$ws.Range("H132").Value = 2162.6128
$ws.Range("I132").Value = 1802.5
$ws.Range("K132").Value = 5407.5
$ws.Range("M132").Value = -2877.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H43").Value = 40000
$ws.Range("J43").Value = 40000
$ws.Range("L43").Value = 40000
$ws.Range("N43").Value = -40386
$ws.Range("H93").Value = 4395.6
$ws.Range("I93").Value = 1598.25
$ws.Range("J93").Value = 6260.5
$ws.Range("K93").Value = 1598.25
$ws.Range("L93").Value = 6260.5
$ws.Range("M93").Value = -350.25
$ws.Range("N93").Value = -8756.5
$ws.Range("H136").Value = 17088.334
$ws.Range("I136").Value = 3816.111
$ws.Range("K136").Value = 11448.333
$ws.Range("M136").Value = -8898.332999999999

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2144.0386
$ws.Range("I132").Value = 1854.238
$ws.Range("K132").Value = 5562.714
$ws.Range("M132").Value = -3032.714
$ws.Range("H136").Value = 4010.3333
$ws.Range("I136").Value = 924.25
$ws.Range("K136").Value = 2772.75
$ws.Range("M136").Value = -222.75
